$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 19 de Julio de 2020 a las 15:41"

# Swap Moldavia (row 62) and Serbia (row 63) country labels.
# Before: row62 = Moldavia, row63 = Serbia
# After:  row62 = Serbia,   row63 = Moldavia
$ws.Cells.Item(62, 1).Value = "Serbia"
$ws.Cells.Item(63, 1).Value = "Moldavia"

# Update numeric statistics for the affected countries/rows.
# Row 4 - Estados Unidos
$ws.Cells.Item(4, 2).Value = 3835430
$ws.Cells.Item(4, 3).Value = 2159
$ws.Cells.Item(4, 5).Value = 1917097
$ws.Cells.Item(4, 7).Value = 6
$ws.Cells.Item(4, 8).Value = 142883

# Row 6 - India
$ws.Cells.Item(6, 2).Value = 1086476
$ws.Cells.Item(6, 3).Value = 8612
$ws.Cells.Item(6, 4).Value = 684662
$ws.Cells.Item(6, 5).Value = 374863
$ws.Cells.Item(6, 7).Value = 123
$ws.Cells.Item(6, 8).Value = 26951

# Row 16 - Arabia Saudita
$ws.Cells.Item(16, 2).Value = 250920
$ws.Cells.Item(16, 3).Value = 2504
$ws.Cells.Item(16, 4).Value = 197735
$ws.Cells.Item(16, 5).Value = 50699
$ws.Cells.Item(16, 7).Value = 39
$ws.Cells.Item(16, 8).Value = 2486

# Row 20 - Alemania
$ws.Cells.Item(20, 2).Value = 202631
$ws.Cells.Item(20, 3).Value = 59
$ws.Cells.Item(20, 5).Value = 5669

# Row 23 - Argentina
$ws.Cells.Item(23, 4).Value = 54105
$ws.Cells.Item(23, 5).Value = 66173
$ws.Cells.Item(23, 7).Value = 26
$ws.Cells.Item(23, 8).Value = 2246

# Row 42 - Paises Bajos
$ws.Cells.Item(42, 2).Value = 51725
$ws.Cells.Item(42, 3).Value = 144

# Row 45 - Portugal
$ws.Cells.Item(45, 2).Value = 48636
$ws.Cells.Item(45, 3).Value = 246
$ws.Cells.Item(45, 4).Value = 33369
$ws.Cells.Item(45, 5).Value = 13578
$ws.Cells.Item(45, 7).Value = 5
$ws.Cells.Item(45, 8).Value = 1689

# Row 62 - now Serbia (new/updated figures)
$ws.Cells.Item(62, 2).Value = 20894
$ws.Cells.Item(62, 3).Value = 396
$ws.Cells.Item(62, 4).Value = 14047
$ws.Cells.Item(62, 5).Value = 6375
$ws.Cells.Item(62, 7).Value = 11
$ws.Cells.Item(62, 8).Value = 472

# Row 63 - now Moldavia (figures that previously belonged to Moldavia/row62)
$ws.Cells.Item(63, 2).Value = 20794
$ws.Cells.Item(63, 3).Value = 0
$ws.Cells.Item(63, 4).Value = 14183
$ws.Cells.Item(63, 5).Value = 5931
$ws.Cells.Item(63, 8).Value = 680

# Row 67 - Uzbekistan
$ws.Cells.Item(67, 2).Value = 16607
$ws.Cells.Item(67, 3).Value = 421
$ws.Cells.Item(67, 4).Value = 9178
$ws.Cells.Item(67, 5).Value = 7344
$ws.Cells.Item(67, 7).Value = 2
$ws.Cells.Item(67, 8).Value = 85

# Row 113 - Sri Lanka
$ws.Cells.Item(113, 2).Value = 2715
$ws.Cells.Item(113, 3).Value = 11
$ws.Cells.Item(113, 5).Value = 669

# Row 141 - Liberia
$ws.Cells.Item(141, 2).Value = 1091
$ws.Cells.Item(141, 3).Value = 3
$ws.Cells.Item(141, 4).Value = 534
$ws.Cells.Item(141, 5).Value = 487

# Row 175 - Islas Feroe
$ws.Cells.Item(175, 2).Value = 191
$ws.Cells.Item(175, 3).Value = 3
$ws.Cells.Item(175, 5).Value = 3
